$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: swap B/C tickers, clear E2
$ws.Range("B2").Value = "NSE:INTELLECT"
$ws.Range("C2").Value = "NSE:AARVI"
$ws.Range("E2").Value = ""

# Row 3: clear B3, replace C3
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = "NSE:CGCL"

# Row 4: clear B4, replace C4
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = "NSE:DPSCLTD"

# Row 5: clear B5, replace C5
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = "NSE:GANGESSECU"

# Row 6
$ws.Range("C6").Value = "NSE:HDFCGROWTH"

# Row 7
$ws.Range("C7").Value = "NSE:ICRA"

# Row 8
$ws.Range("C8").Value = "NSE:JHS"

# Row 9
$ws.Range("C9").Value = "NSE:NITIRAJ"

# Row 10
$ws.Range("C10").Value = "NSE:ONMOBILE"

# Row 11
$ws.Range("C11").Value = "NSE:PVP"

# Row 12
$ws.Range("C12").Value = "NSE:SAIL"

# Remove old rows 13-22 entirely (they no longer exist in the new data)
$ws.Rows("13:22").Delete()
